$d = $word.ActiveDocument

# Replace the ${kelurahan} placeholder with the literal value "Leuwigajah",
# mirroring the Word edit that hard-coded the kelurahan name for this
# particular document instance.
$target = $d.Content
$replaced = $target.Find.Execute("`${kelurahan}", $true, $false, $false, $false, $false,
                                  $true, 1, $false, "Leuwigajah", 2)

# After an in-place edit, Word repositions the automatic "_GoBack" bookmark
# to mark the last editing location -- a collapsed bookmark sitting right
# after the text that was just typed/replaced.
$edited = $d.Content
$found = $edited.Find.Execute("Leuwigajah", $true, $false, $false, $false, $false,
                               $true, 1, $false, "", 0)
if ($found) {
    $goBack = $d.Range($edited.End, $edited.End)
    $d.Bookmarks.Add("_GoBack", $goBack) | Out-Null
}

Write-Output ("kelurahan replaced: " + $replaced + "; _GoBack bookmark set: " + $found)
